# Normalise the resolution strings in column D ("resolusi_layar") from
# "1920 , 1080" style formatting to "1920x1080" style formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, fix item 4 (row 4) by hand.
$cell4 = $ws.Cells.Item(4, 4)
$cell4.Value2 = $cell4.Value2.Replace(" , ", "x")

# Then replace the rest of the column in one pass.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 4) { continue }
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v.Replace(" , ", "x")
    }
}

$ws.Range("D6").Select()
